$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear A5:B9 (these columns have no data in rows 5-9 after the edit)
$ws.Range("A5:B9").ClearContents()

# Row 1 - header
$ws.Range("A1").Value = "labelA_stimuli"
$ws.Range("B1").Value = "labelB_stimuli"
$ws.Range("C1").Value = "targetA_stimuli"
$ws.Range("D1").Value = "targetB_stimuli"
$ws.Range("E1").Value = "labelA_image_stimuli"
$ws.Range("F1").Value = "labelB_image_stimuli"
$ws.Range("G1").Value = "targetA_image_stimuli"
$ws.Range("H1").Value = "targetB_image_stimuli"

# Rows 2-4 (A,B still " ")
$ws.Range("A2").Value = " "
$ws.Range("B2").Value = " "
$ws.Range("C2").Value = "happy"
$ws.Range("D2").Value = "pain"
$ws.Range("E2").Value = "flower1.jpg"
$ws.Range("F2").Value = "insect1.jpg"
$ws.Range("G2").Value = "blank.jpg"
$ws.Range("H2").Value = "blank.jpg"
$ws.Range("A3").Value = " "
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = "friendly"
$ws.Range("D3").Value = "mean"
$ws.Range("E3").Value = "flower2.jpg"
$ws.Range("F3").Value = "insect2.jpg"
$ws.Range("G3").Value = "blank.jpg"
$ws.Range("H3").Value = "blank.jpg"
$ws.Range("A4").Value = " "
$ws.Range("B4").Value = " "
$ws.Range("C4").Value = "pretty"
$ws.Range("D4").Value = "hostile"
$ws.Range("E4").Value = "flower3.jpg"
$ws.Range("F4").Value = "insect3.jpg"
$ws.Range("G4").Value = "blank.jpg"
$ws.Range("H4").Value = "blank.jpg"

# Rows 5-9 (C,D,E,F,G,H only; A,B cleared above)
$ws.Range("C5").Value = "good"
$ws.Range("D5").Value = "hateful"
$ws.Range("E5").Value = "flower4.jpg"
$ws.Range("F5").Value = "insect4.jpg"
$ws.Range("G5").Value = "blank.jpg"
$ws.Range("H5").Value = "blank.jpg"
$ws.Range("C6").Value = "happy"
$ws.Range("D6").Value = "pain"
$ws.Range("E6").Value = "flower1.jpg"
$ws.Range("F6").Value = "insect1.jpg"
$ws.Range("G6").Value = "blank.jpg"
$ws.Range("H6").Value = "blank.jpg"
$ws.Range("C7").Value = "friendly"
$ws.Range("D7").Value = "mean"
$ws.Range("E7").Value = "flower2.jpg"
$ws.Range("F7").Value = "insect2.jpg"
$ws.Range("G7").Value = "blank.jpg"
$ws.Range("H7").Value = "blank.jpg"
$ws.Range("C8").Value = "pretty"
$ws.Range("D8").Value = "hostile"
$ws.Range("E8").Value = "flower3.jpg"
$ws.Range("F8").Value = "insect3.jpg"
$ws.Range("G8").Value = "blank.jpg"
$ws.Range("H8").Value = "blank.jpg"
$ws.Range("C9").Value = "good"
$ws.Range("D9").Value = "hateful"
$ws.Range("E9").Value = "flower4.jpg"
$ws.Range("F9").Value = "insect4.jpg"
$ws.Range("G9").Value = "blank.jpg"
$ws.Range("H9").Value = "blank.jpg"

# Column D (rows 2-9) gets Text number format (numFmtId 49)
$ws.Range("D2:D9").NumberFormat = "@"

# Update selection to match the target view state
$ws.Range("E13").Select()
